# "added more data points"
# Adds three new rows (45-47) of travel data to Sheet1:
#   45: USA / CA / Los Angeles / (last visited 2022-03-01) / comment
#   46: New Zealand / Bay of Plenty / Rotorua / (last visited 2020-03-01)
#   47: New Zealand / Tasman / Nelson / (last visited 2015-12-01)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45: USA / CA / Los Angeles -----------------------------------
$ws.Range("A45").Value = "USA"
$ws.Range("B45").Value = "CA "
$ws.Range("C45").Value = "Los Angeles"
$ws.Range("E45").NumberFormat = $ws.Range("E44").NumberFormat
$ws.Range("E45").Value = (Get-Date -Year 2022 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F45").Value = 'I really only transit through LAX, I don''t think I''ve ever actally stepped foot in LA…'

# --- Row 46: New Zealand / Bay of Plenty / Rotorua ---------------------
$ws.Range("A46").Value = "New Zealand"
$ws.Range("C46").Value = "Rotorua"
$ws.Range("B46").Value = "Bay of Plenty"
$ws.Range("E46").NumberFormat = $ws.Range("E44").NumberFormat
$ws.Range("E46").Value = (Get-Date -Year 2020 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0)

# --- Row 47: New Zealand / Tasman / Nelson ------------------------------
$ws.Range("A47").Value = "New Zealand"
$ws.Range("C47").Value = "Nelson"
$ws.Range("B47").Value = "Tasman"
$ws.Range("E47").NumberFormat = $ws.Range("E44").NumberFormat
$ws.Range("E47").Value = (Get-Date -Year 2015 -Month 12 -Day 1 -Hour 0 -Minute 0 -Second 0)

# --- View state: scroll down a bit and move the active selection -------
$win = $excel.ActiveWindow
$win.Left = 120
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("B48").Select() | Out-Null
